# Apply updates to column F (dSF) values as per commit:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F4").Value = -1
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("F17").Value = 2
$ws.Range("F20").Value = -2
